$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 becomes a new "Bigmac" line item -----------------------------
$ws.Range("B6").Value = "Bigmac"
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 41

# --- Change the profit-margin formula (3.4% cut instead of *0.966) ------
# E2 carries its own (non-shared) formula.
$ws.Range("E2").Formula = "=D2-(D2*3.4%)-2.6-C2"

# E3:E13 is one shared-formula block (master at E3); writing the whole
# range in one call keeps it a single shared formula instead of exploding
# it into per-cell formulas.
$ws.Range("E3:E13").Formula = "=D3-(D3*3.4%)-2.6-C3"

# --- Move the active selection -------------------------------------------
$ws.Range("D19").Select()

# --- Misc window / file metadata (as recorded by Excel on save) ---------
$excel.Left = 2500
$excel.Top = 3020
$wb.Windows.Item(1).Left = 2500
$wb.Windows.Item(1).Top = 3020
$wb.Path = "/Users/maxy-macbook/Desktop/gp-food.github.io-master/"
